$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# --- Row 17: Brooklyn Nets @ Los Angeles Clippers (final score + OT flag) ---
$ws.Range("D17").Value = 114
$ws.Range("F17").Value = 125
$ws.Range("G17").Value = "NA"

# --- Row 18: Miami Heat @ Orlando Magic ---
$ws.Range("D18").Value = 87
$ws.Range("F18").Value = 105

# --- Row 19: Denver Nuggets @ Washington Wizards ---
$ws.Range("D19").Value = 113
$ws.Range("F19").Value = 104

# --- Row 20: Boston Celtics @ Houston Rockets ---
$ws.Range("D20").Value = 116
$ws.Range("F20").Value = 107

# --- Row 21: Indiana Pacers @ Phoenix Suns ---
$ws.Range("D21").Value = 110
$ws.Range("F21").Value = 117

# --- Row 22: Portland Trail Blazers @ Los Angeles Lakers ---
$ws.Range("D22").Value = 110
$ws.Range("F22").Value = 134

# --- Row 23: Milwaukee Bucks @ Detroit Pistons (score + forecast outcome) ---
$ws.Range("D23").Value = 122
$ws.Range("F23").Value = 113
$ws.Range("I23").Value = "Milwaukee Bucks"
$ws.Range("J23").Value = "Detroit Pistons"
$ws.Range("K23").Value = "Detroit Pistons"
$ws.Range("L23").Value = "No"

# --- Row 24: Cleveland Cavaliers @ Orlando Magic ---
$ws.Range("D24").Value = 126
$ws.Range("F24").Value = 99
$ws.Range("I24").Value = "Cleveland Cavaliers"
$ws.Range("J24").Value = "Orlando Magic"
$ws.Range("K24").Value = "Cleveland Cavaliers"
$ws.Range("L24").Value = "Yes"

# --- Row 25: San Antonio Spurs @ Philadelphia 76ers ---
$ws.Range("D25").Value = 123
$ws.Range("F25").Value = 133
$ws.Range("I25").Value = "Philadelphia 76ers"
$ws.Range("J25").Value = "San Antonio Spurs"
$ws.Range("K25").Value = "San Antonio Spurs"
$ws.Range("L25").Value = "No"

# --- Row 26: Memphis Grizzlies @ Toronto Raptors ---
$ws.Range("D26").Value = 108
$ws.Range("F26").Value = 100
$ws.Range("I26").Value = "Memphis Grizzlies"
$ws.Range("J26").Value = "Toronto Raptors"
$ws.Range("K26").Value = "Toronto Raptors"
$ws.Range("L26").Value = "No"

# --- Row 27: Charlotte Hornets @ Minnesota Timberwolves ---
$ws.Range("D27").Value = 128
$ws.Range("F27").Value = 125
$ws.Range("I27").Value = "Charlotte Hornets"
$ws.Range("J27").Value = "Minnesota Timberwolves"
$ws.Range("K27").Value = "Charlotte Hornets"
$ws.Range("L27").Value = "Yes"

# --- Row 28: Boston Celtics @ Dallas Mavericks ---
$ws.Range("D28").Value = 119
$ws.Range("F28").Value = 110
$ws.Range("I28").Value = "Boston Celtics"
$ws.Range("J28").Value = "Dallas Mavericks"
$ws.Range("K28").Value = "Boston Celtics"
$ws.Range("L28").Value = "Yes"

# --- Row 29: Chicago Bulls @ Phoenix Suns ---
$ws.Range("D29").Value = 113
$ws.Range("F29").Value = 115
$ws.Range("I29").Value = "Phoenix Suns"
$ws.Range("J29").Value = "Chicago Bulls"
$ws.Range("K29").Value = "Chicago Bulls"
$ws.Range("L29").Value = "No"

# --- Row 30: Sacramento Kings @ Atlanta Hawks ---
$ws.Range("D30").Value = 107
$ws.Range("F30").Value = 122
$ws.Range("I30").Value = "Sacramento Kings"
$ws.Range("J30").Value = "Atlanta Hawks"
$ws.Range("K30").Value = "Atlanta Hawks"
$ws.Range("L30").Value = "No"

# Highlight the correctly-forecasted winner in column K (matches the green
# "correct pick" styling already used elsewhere in the sheet, e.g. K20/K21)
$ws.Range("K20").Copy()
$ws.Range("K24").PasteSpecial(-4122)
$ws.Range("K27").PasteSpecial(-4122)
$ws.Range("K28").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Scroll the view down and leave H18 as the active selected cell
$ws.Range("H18").Select()
